$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel; force them to remain text
# by temporarily switching to Text format, then clearing the style again
# so no stray style index is left attached to the cell.
$numericLookingCells = @("D4", "D5", "D6", "D9", "D10", "D11", "D15", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D43", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated "Price" column (D) values
$ws.Range("D2").Value = "26.086.52"
$ws.Range("D3").Value = "1.666.75"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").Value = "209.86"
$ws.Range("D6").Value = "0.5259"
$ws.Range("D9").Value = "0.06286"
$ws.Range("D10").Value = "21.17"
$ws.Range("D11").Value = "0.07528"
$ws.Range("D12").Value = "1.680.40"
$ws.Range("D15").Value = "66.60"
$ws.Range("D17").Value = "26.139.95"
$ws.Range("D19").Value = "4.717"
$ws.Range("D20").Value = "186.08"
$ws.Range("D22").Value = "6.154"
$ws.Range("D24").Value = "149.82"
$ws.Range("D25").Value = "0.1245"
$ws.Range("D26").Value = "7.472"
$ws.Range("D27").Value = "15.89"
$ws.Range("D28").Value = "0.06293"
$ws.Range("D30").Value = "1.275"
$ws.Range("D31").Value = "3.494"
$ws.Range("D32").Value = "3.407"
$ws.Range("D33").Value = "1.633"
$ws.Range("D34").Value = "0.9955"
$ws.Range("D35").Value = "0.6032"
$ws.Range("D36").Value = "2.408"
$ws.Range("D37").Value = "2.726"
$ws.Range("D38").Value = "1.107.81"
$ws.Range("D40").Value = "0.01614"
$ws.Range("D41").Value = "0.8746"
$ws.Range("D43").Value = "99.76"
$ws.Range("D46").Value = "55.37"
$ws.Range("D47").Value = "0.9997"
$ws.Range("D48").Value = "8.034"
$ws.Range("D50").Value = "0.4244"
$ws.Range("D51").Value = "5.966"

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}

# Updated "Volume(1h)" column (E) values
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E8").Value = "  -3.24%  "
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("E14").Value = "  -4.10%  "
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("E16").Value = "  -5.10%  "
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  -4.75%  "
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E51").Value = "  -1.42%  "
